$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-15) were re-shuffled: each row's Fecha (D), Volumen (M),
# Precio minimo (N), Precio maximo (O), Precio promedio ponderado (P),
# Origen (R) and Precio $/Kg (S) are updated to reflect the new weekly
# arrangement of records. Columns A,B,C,E,F,G,H,I,J,K,L,Q,T are unchanged.

$rowData = @{
    2  = @{ D = 44208; M = 85;  N = 3000; O = 3000; P = 3000; R = "Provincia de Linares"; S = 1500 }
    4  = @{ D = 44188; M = $null; N = 3000; O = $null; P = 3240; R = "Provincia de Linares"; S = 1620 }
    5  = @{ D = 44617; M = 90;  N = 6500; O = 6500; P = 6500; R = "Provincia de Curicó"; S = 3250 }
    6  = @{ D = 44533; M = 150; N = 4000; O = 4000; P = 4000; R = "Provincia de Curicó"; S = 2000 }
    7  = @{ D = 44586; M = 250; N = 5000; O = 5000; P = 5000; R = "Provincia de Curicó"; S = 2500 }
    8  = @{ D = 44194; M = 120; N = 3000; O = 3000; P = 3000; R = "Provincia de Linares"; S = 1500 }
    9  = @{ D = 44231; M = 150; N = 3400; O = 3400; P = 3400; R = $null; S = 1700 }
    10 = @{ D = 44236; M = 300; N = 3600; O = 4000; P = 3800; R = $null; S = 1900 }
    11 = @{ D = 44238; M = 300; N = 3600; O = 4000; P = 3800; R = $null; S = 1900 }
    12 = @{ D = 44174; M = 200; N = 3200; O = 3200; P = 3200; R = "Provincia de Curicó"; S = 1600 }
    14 = @{ D = 44232; M = 200; N = 3000; O = 3000; P = 3000; R = $null; S = 1500 }
    15 = @{ D = 44168; M = 170; N = 8000; O = 8000; P = 8000; R = "Provincia de Linares"; S = 4000 }
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]

    $ws.Cells.Item($row, 4).Value = $vals.D

    if ($vals.M -ne $null) { $ws.Cells.Item($row, 13).Value = $vals.M }
    $ws.Cells.Item($row, 14).Value = $vals.N
    if ($vals.O -ne $null) { $ws.Cells.Item($row, 15).Value = $vals.O }
    $ws.Cells.Item($row, 16).Value = $vals.P
    if ($vals.R -ne $null) { $ws.Cells.Item($row, 18).Value = $vals.R }
    $ws.Cells.Item($row, 19).Value = $vals.S
}
